$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1958762886597938
$ws.Range("C2").Value = 0.5532646048109966
$ws.Range("J2").Value = 0.0274914089347079
$ws.Range("P2").Value = 0.1305841924398626
$ws.Range("S2").Value = 0.09278350515463918
$ws.Range("B3").Value = 0.02247191011235955
$ws.Range("C3").Value = 0.07303370786516854
$ws.Range("J3").Value = 0.01685393258426966
$ws.Range("P3").Value = 0.7471910112359551
$ws.Range("S3").Value = 0.1404494382022472
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.05150214592274678
$ws.Range("D6").Value = 0.008583690987124463
$ws.Range("F6").Value = 0.07296137339055794
$ws.Range("J6").Value = 0.2532188841201717
$ws.Range("O6").Value = 0.03004291845493562
$ws.Range("Q6").Value = 0.1759656652360515
$ws.Range("R6").Value = 0.0815450643776824
$ws.Range("S6").Value = 0.3261802575107296
$ws.Range("B7").Value = 0.08658008658008658
$ws.Range("D7").Value = 0.02164502164502164
$ws.Range("E7").Value = 0.004329004329004329
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.1558441558441558
$ws.Range("O7").Value = 0.02597402597402598
$ws.Range("Q7").Value = 0.1991341991341991
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.3679653679653679
$ws.Range("B8").Value = 0.07407407407407407
$ws.Range("D8").Value = 0.01296296296296296
$ws.Range("E8").Value = 0.001851851851851852
$ws.Range("F8").Value = 0.05555555555555555
$ws.Range("J8").Value = 0.15
$ws.Range("O8").Value = 0.02407407407407407
$ws.Range("Q8").Value = 0.1462962962962963
$ws.Range("R8").Value = 0.09444444444444444
$ws.Range("S8").Value = 0.4407407407407408
$ws.Range("B9").Value = 0.09844559585492228
$ws.Range("D9").Value = 0.0155440414507772
$ws.Range("F9").Value = 0.07772020725388601
$ws.Range("J9").Value = 0.1398963730569948
$ws.Range("O9").Value = 0.02072538860103627
$ws.Range("Q9").Value = 0.1450777202072539
$ws.Range("R9").Value = 0.1347150259067358
$ws.Range("S9").Value = 0.3678756476683938
$ws.Range("B10").Value = 0.09777138749101366
$ws.Range("D10").Value = 0.01869158878504673
$ws.Range("E10").Value = 0.001437814521926672
$ws.Range("F10").Value = 0.06326383896477354
$ws.Range("J10").Value = 0.1301222142343638
$ws.Range("O10").Value = 0.01437814521926671
$ws.Range("Q10").Value = 0.205607476635514
$ws.Range("R10").Value = 0.08986340762041696
$ws.Range("S10").Value = 0.3788641265276779
$ws.Range("G11").Value = 0.1421052631578947
$ws.Range("J11").Value = 0.08947368421052632
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5526315789473685
$ws.Range("S11").Value = 0.01578947368421053
$ws.Range("G12").Value = 0.7534883720930232
$ws.Range("J12").Value = 0.1581395348837209
$ws.Range("K12").Value = 0.009302325581395349
$ws.Range("L12").Value = 0.03720930232558139
$ws.Range("S12").Value = 0.04186046511627907
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2571428571428571
$ws.Range("S13").Value = 0.02857142857142857
$ws.Range("F15").Value = 0.0326530612244898
$ws.Range("H15").Value = 0.1918367346938775
$ws.Range("I15").Value = 0.06122448979591837
$ws.Range("J15").Value = 0.3265306122448979
$ws.Range("K15").Value = 0.08979591836734693
$ws.Range("M15").Value = 0.01224489795918367
$ws.Range("N15").Value = 0.004081632653061225
$ws.Range("O15").Value = 0.0653061224489796
$ws.Range("S15").Value = 0.2163265306122449
$ws.Range("F16").Value = 0.02604166666666667
$ws.Range("H16").Value = 0.1979166666666667
$ws.Range("I16").Value = 0.04166666666666666
$ws.Range("J16").Value = 0.4010416666666667
$ws.Range("K16").Value = 0.1145833333333333
$ws.Range("M16").Value = 0.03645833333333334
$ws.Range("N16").Value = 0.005208333333333333
$ws.Range("O16").Value = 0.046875
$ws.Range("S16").Value = 0.1302083333333333
$ws.Range("F17").Value = 0.01461377870563674
$ws.Range("H17").Value = 0.2025052192066806
$ws.Range("I17").Value = 0.09603340292275574
$ws.Range("J17").Value = 0.4217118997912317
$ws.Range("K17").Value = 0.0918580375782881
$ws.Range("M17").Value = 0.01670146137787056
$ws.Range("O17").Value = 0.05845511482254697
$ws.Range("S17").Value = 0.09812108559498957
$ws.Range("F18").Value = 0.01260504201680672
$ws.Range("H18").Value = 0.2184873949579832
$ws.Range("I18").Value = 0.08823529411764706
$ws.Range("J18").Value = 0.3781512605042017
$ws.Range("K18").Value = 0.1050420168067227
$ws.Range("M18").Value = 0.02100840336134454
$ws.Range("N18").Value = 0.004201680672268907
$ws.Range("O18").Value = 0.05042016806722689
$ws.Range("S18").Value = 0.1218487394957983
$ws.Range("F19").Value = 0.01529497450837582
$ws.Range("H19").Value = 0.2272396212672979
$ws.Range("I19").Value = 0.0764748725418791
$ws.Range("J19").Value = 0.3525127458120903
$ws.Range("K19").Value = 0.1318281136198106
$ws.Range("M19").Value = 0.01092498179169701
$ws.Range("N19").Value = 0.0007283321194464676
$ws.Range("O19").Value = 0.07283321194464676
$ws.Range("S19").Value = 0.112163146394756
